$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57 (shifts existing rows 57-70 down to 58-71),
# mirroring the copied-down formatting from the row above (keeps the
# date style on column D).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Range("A57").Value = 1
$ws.Range("B57").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C57").Value = 'Arica y Parinacota'
$ws.Range("D57").Value = 44889
$ws.Range("E57").Value = 15
$ws.Range("F57").Value = 100112031
$ws.Range("G57").Value = 'Poroto verde'
$ws.Range("H57").Value = 'Sin especificar'
$ws.Range("I57").Value = 'Primera'
$ws.Range("J57").Value = 500
$ws.Range("K57").Value = 900
$ws.Range("L57").Value = 1000
$ws.Range("M57").Value = 950
$ws.Range("N57").Value = '$/kilo'
$ws.Range("O57").Value = 'Región de Arica y Parinacota'
$ws.Range("P57").Value = 950
$ws.Range("Q57").Value = 1
$ws.Range("R57").Value = 'Hortaliza'
